$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values change to what used to be Row 4's values
$ws.Range("D2").Value = 44518
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 667

# Row 4 values change to what used to be Row 2's values
$ws.Range("D4").Value = 44525
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 533
